# Updates cryptos list (Price / Volume(1h) columns) as scraped on
# Mon Feb 27 16:55:01 UTC 2023 with GitHub Actions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row=2;  D="23.468.15";   E="  +1.07%  " },
    @{ Row=3;  D="1.643.03";    E="  +2.40%  " },
    @{ Row=4;  D=$null;         E="  +0.25%  " },
    @{ Row=5;  D=$null;         E="  +0.26%  " },
    @{ Row=6;  D="303.89";      E="  -0.26%  " },
    @{ Row=7;  D="0.3769";      E="  -0.06%  " },
    @{ Row=8;  D="52.19";       E="  -0.95%  " },
    @{ Row=9;  D="0.3651";      E="  +0.70%  " },
    @{ Row=10; D=$null;         E="  -1.43%  " },
    @{ Row=11; D="0.08141";     E="  -0.05%  " },
    @{ Row=12; D=$null;         E="  +0.24%  " },
    @{ Row=13; D=$null;         E="  +0.39%  " },
    @{ Row=14; D="6.664";       E=$null },
    @{ Row=15; D="0.00001262";  E="  +1.26%  " },
    @{ Row=16; D="7.313";       E="  -0.61%  " },
    @{ Row=17; D="1.642.16";    E="  +2.38%  " },
    @{ Row=18; D="94.39";       E="  +0.34%  " },
    @{ Row=19; D="0.06936";     E="  -0.02%  " },
    @{ Row=20; D="18.20";       E="  +0.39%  " },
    @{ Row=21; D="6.571";       E="  +0.58%  " },
    @{ Row=22; D="1.005";       E="  +0.25%  " },
    @{ Row=23; D="23.478.90";   E="  +1.09%  " },
    @{ Row=24; D="12.90";       E=$null },
    @{ Row=25; D="3.260";       E="  +6.46%  " },
    @{ Row=26; D="2.444";       E="  +0.52%  " },
    @{ Row=27; D="21.27";       E="  +0.42%  " },
    @{ Row=28; D="151.37";      E="  +0.57%  " },
    @{ Row=29; D="5.323";       E="  +0.56%  " },
    @{ Row=30; D="136.39";      E="  +1.11%  " },
    @{ Row=31; D=$null;         E="  -3.64%  " },
    @{ Row=32; D="1.826.90";    E="  +2.75%  " },
    @{ Row=33; D="6.937";       E="  +2.54%  " },
    @{ Row=34; D="11.02";       E="  +6.75%  " },
    @{ Row=35; D="0.9664";      E="  +1.09%  " },
    @{ Row=36; D="0.02874";     E="  +3.94%  " },
    @{ Row=37; D="6.295";       E="  +2.88%  " },
    @{ Row=38; D="0.2570";      E="  +2.13%  " },
    @{ Row=39; D="0.07278";     E="  -2.18%  " },
    @{ Row=40; D="0.08845";     E="  +0.81%  " },
    @{ Row=41; D="1.381";       E="  -2.00%  " },
    @{ Row=42; D="0.7154";      E="  +0.76%  " },
    @{ Row=43; D="16.47";       E="  +3.91%  " },
    @{ Row=44; D="12.62";       E="  +1.55%  " },
    @{ Row=45; D="0.6587";      E="  +0.86%  " },
    @{ Row=46; D=$null;         E="  +1.68%  " },
    @{ Row=47; D="1.003";       E="  +0.33%  " },
    @{ Row=48; D="3.998";       E="  -0.31%  " },
    @{ Row=49; D="0.08015";     E="  +0.82%  " },
    @{ Row=50; D="1.219";       E="  +1.62%  " },
    @{ Row=51; D="128.32";      E="  -4.14%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Column D ("Price") holds text (not numeric) values in the source
        # data, e.g. "23.468.15" or "303.89". Force the cell to remain
        # text so Excel doesn't reinterpret plain-looking numbers (like
        # "303.89") as actual numeric values, then restore the default
        # "Normal" style so no stray formatting is introduced.
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
